$d = $word.ActiveDocument

# The second paragraph holds a Word field whose code spells out
#   { m:'doc.html'.fromHTMLURI() }
# via fldChar/instrText runs. The commit replaces that field with plain
# literal text runs spelling out the same characters (no field anymore),
# keeping the _GoBack bookmark in place around "doc.html".
#
# Field codes/instrText are not part of Range.Text (they carry no visible
# text), so Find/Range.Text based editing cannot target them precisely.
# Range.InsertXML, however, replaces the *entire* contents of the Range it
# is called on with the OOXML we provide, which lets us swap the whole
# paragraph's run sequence in one precise operation.

# Locate the paragraph that actually contains the field (rather than
# hardcoding its index) by checking Range.Fields.Count per paragraph -
# Field objects derived via $d.Fields.Item(n).Code/.Paragraphs resolve to
# the wrong Range in this environment, so we avoid that path entirely.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Fields.Count -gt 0) {
        $target = $candidate.Range
        break
    }
}
if ($null -eq $target) {
    throw "Could not locate the paragraph containing the m:'doc.html' field"
}

$newParaXml = '<w:p ' +
  'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" ' +
  'w:rsidR="00C52979" w:rsidRDefault="00C52979" w:rsidP="00F5495F">' +
    '<w:r><w:t>{</w:t></w:r>' +
    '<w:r><w:t>m</w:t></w:r>' +
    '<w:r><w:t>:</w:t></w:r>' +
    '<w:r><w:t>''</w:t></w:r>' +
    '<w:r><w:t>doc.html</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
    '<w:bookmarkEnd w:id="0"/>' +
    '<w:r><w:t>''.fromHTMLURI()</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">}</w:t></w:r>' +
  '</w:p>'

$target.InsertXML($newParaXml)
